$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = 12
    3  = -1
    4  = -2
    5  = 2
    6  = -5
    8  = -1
    9  = 1
    10 = 3
    12 = 1
    13 = -1
    14 = -1
    15 = 1
    17 = -1
    18 = 2
    19 = 1
    22 = 1
    23 = 1
    24 = -1
    25 = -7
    26 = -1
    27 = -2
    28 = 3
    29 = -3
    30 = -1
}

foreach ($row in $newValues.Keys) {
    $ws.Range("F$row").Value = $newValues[$row]
}
